$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the week's dates forward by 7 days (new week: Oct 17 - Oct 23)
$ws.Range("B5").Value = 44851
$ws.Range("C5").Value = 44852
$ws.Range("D5").Value = 44853
$ws.Range("E5").Value = 44854
$ws.Range("F5").Value = 44855
$ws.Range("G5").Value = 44856
$ws.Range("H5").Value = 44857

# Clear out hours logged in the previous week that no longer apply
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("F9").ClearContents()

# Update TA Meeting hours on Friday (H11) from 0.25 to 1
$ws.Range("H11").Value = 1

# Update the active selection to reflect where the user left off
$ws.Range("D7").Select()
